$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "hh"
$ws.Range("C2").Value = "h"
